$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 new data (entered first so new shared strings are interned in
#     the same order as the target: First2D..., Added..., Mel conv2d...) ---
$ws.Range("A3").Value = "First 2D Convolutional Attempt on Single Tagged Clips"
$ws.Range("B3").Value = "Added a 2D Conv and Leaky Relu"
$ws.Range("C3").Value = "The 444 clips of 10 seconds with one tag"
$ws.Range("D3").Value = "Mel, conv2d, leakyrelu, flatten, dense"
$ws.Range("E3").Value = 3.6
$ws.Range("F3").Value = 0.77
$ws.Range("G3").Value = 48
$ws.Range("H3").Value = 11
$ws.Range("I3").Value = 14
$ws.Range("J3").Value = 38

$ws.Range("K3").Formula = "=(G3+J3)/SUM(G3:J3)"
$ws.Range("K3").NumberFormat = "0.00"

# --- Row 1 header changes ---
# G1 header: "True Positives" -> "True Negatives" (style unchanged, s="2")
$ws.Range("G1").Value = "True Negatives"

# K1 new header: "Actual Accuracy", matching style of A1 (s="1", bold)
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Actual Accuracy"

# --- Row 2 updates ---
$ws.Range("G2").Value = 44
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 31

$ws.Range("K2").Formula = "=(G2+J2)/SUM(G2:J2)"
$ws.Range("K2").NumberFormat = "0.00"

# --- Column B width ---
# The runtime's "characters -> stored width" conversion adds ~5/6 (0.8333)
# to whatever ColumnWidth we set (stored = chars + 5/MDW with MDW=6), so to
# land on the target stored width of exactly 23 we back it out here.
$ws.Columns.Item(2).ColumnWidth = 22.17

# --- Selection ---
$ws.Range("H3").Select() | Out-Null
